# "wrapping up test file audit"
# Remove the stray "Sheet" parameter row (row 16: Sheet / 3 / 4) from the
# optimization_parameters worksheet; the row below it (simulation_timepoints)
# shifts up to take its place.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")
$ws.Rows.Item(16).Delete()
$ws.Rows.Item(16).Select()

# Final active sheet/tab ends up on the last sheet (optimization_diagnostics).
$lastSheet = $wb.Worksheets.Item("optimization_diagnostics")
$lastSheet.Activate()
